$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 43-48 (no longer present in the updated dataset)
$ws.Range("A43:A48").EntireRow.Delete()

# Update rows 2-42 with the refreshed article data
$ws.Range("A2").Value = '4813022'
$ws.Range("B2").Value = 'https://vnexpress.net/nga-hoai-nghi-kha-nang-ong-trump-giai-quyet-xung-dot-ukraine-4813022.html'
$ws.Range("C2").Value = 'Nga hoài nghi khả năng ông Trump giải quyết xung đột Ukraine'
$ws.Range("A3").Value = '4812992'
$ws.Range("B3").Value = 'https://vnexpress.net/chu-tich-nuoc-luong-cuong-sap-du-apec-tai-peru-4812992.html'
$ws.Range("C3").Value = 'Chủ tịch nước Lương Cường sắp dự APEC tại Peru'
$ws.Range("A4").Value = '4812836'
$ws.Range("B4").Value = 'https://vnexpress.net/nguoi-ung-ho-trump-reo-ho-bat-khoc-an-mung-4812836.html'
$ws.Range("C4").Value = 'Người ủng hộ Trump reo hò, bật khóc ăn mừng'
$ws.Range("A5").Value = '4812870'
$ws.Range("B5").Value = 'https://vnexpress.net/lanh-dao-the-gioi-chuc-mung-ong-trump-4812870.html'
$ws.Range("C5").Value = 'Lãnh đạo thế giới chúc mừng ông Trump'
$ws.Range("A6").Value = '4812802'
$ws.Range("B6").Value = 'https://vnexpress.net/nguoi-tre-viet-hao-hung-theo-doi-bau-cu-tong-thong-my-4812802.html'
$ws.Range("C6").Value = 'Người trẻ Việt hào hứng theo dõi bầu cử tổng thống Mỹ'
$ws.Range("A7").Value = '4812787'
$ws.Range("B7").Value = 'https://vnexpress.net/dang-cong-hoa-gianh-da-so-thuong-vien-my-4812787.html'
$ws.Range("C7").Value = 'Đảng Cộng hòa giành đa số Thượng viện Mỹ'
$ws.Range("A8").Value = '4812603'
$ws.Range("B8").Value = 'https://vnexpress.net/cu-tri-kien-elon-musk-vi-khong-chon-ngau-nhien-nguoi-nhan-mot-trieu-usd-4812603.html'
$ws.Range("C8").Value = 'Cử tri kiện Elon Musk vì không chọn ngẫu nhiên người nhận một triệu USD'
$ws.Range("A9").Value = '4810598'
$ws.Range("B9").Value = 'https://vnexpress.net/giang-co-o-7-bang-chien-truong-dinh-doat-bau-cu-my-4810598.html'
$ws.Range("C9").Value = 'Ông Trump áp đảo ở các bang chiến trường'
$ws.Range("A10").Value = '4812723'
$ws.Range("B10").Value = 'https://vnexpress.net/tam-trang-nguoi-my-trong-luc-cho-ket-qua-bau-cu-4812723.html'
$ws.Range("C10").Value = 'Tâm trạng người Mỹ trong lúc chờ kết quả bầu cử'
$ws.Range("A11").Value = '4812730'
$ws.Range("B11").Value = 'https://vnexpress.net/barron-trump-lan-dau-di-bo-phieu-4812730.html'
$ws.Range("C11").Value = 'Barron Trump lần đầu đi bỏ phiếu'
$ws.Range("A12").Value = '4812674'
$ws.Range("B12").Value = 'https://vnexpress.net/vui-buon-cua-phe-cong-hoa-dan-chu-khi-ong-trump-dan-truoc-4812674.html'
$ws.Range("C12").Value = 'Vui, buồn của phe Cộng hòa - Dân chủ khi ông Trump dẫn trước'
$ws.Range("A13").Value = '4812752'
$ws.Range("B13").Value = 'https://vnexpress.net/nevada-gap-rac-roi-vi-cu-tri-tre-khong-biet-ky-ten-4812752.html'
$ws.Range("C13").Value = 'Nevada gặp rắc rối vì cử tri trẻ không biết ký tên'
$ws.Range("A14").Value = '4812592'
$ws.Range("B14").Value = 'https://vnexpress.net/nuoc-my-trong-ngay-bau-cu-4812592.html'
$ws.Range("C14").Value = 'Nước Mỹ trong ngày bầu cử'
$ws.Range("A15").Value = '4812571'
$ws.Range("B15").Value = 'https://vnexpress.net/nhung-nguoi-nhap-cu-my-neu-ly-do-bau-cho-ong-trump-4812571.html'
$ws.Range("C15").Value = 'Những người nhập cư Mỹ nêu lý do bầu cho ông Trump'
$ws.Range("A16").Value = '4812584'
$ws.Range("B16").Value = 'https://vnexpress.net/cu-tri-ung-ho-ong-trump-ve-kinh-te-ba-harris-ve-quyen-pha-thai-4812584.html'
$ws.Range("C16").Value = 'Cử tri ủng hộ ông Trump về kinh tế, bà Harris về quyền phá thai'
$ws.Range("A17").Value = '4811304'
$ws.Range("B17").Value = 'https://vnexpress.net/truyen-thong-my-xuong-ten-nguoi-dac-cu-the-nao-4811304.html'
$ws.Range("C17").Value = 'Truyền thông Mỹ xướng tên người đắc cử thế nào'
$ws.Range("A18").Value = '4812306'
$ws.Range("B18").Value = 'https://vnexpress.net/nguoi-my-goc-viet-neu-ly-do-bau-cho-trump-harris-4812306.html'
$ws.Range("C18").Value = 'Người Mỹ gốc Việt nêu lý do bầu cho Trump, Harris'
$ws.Range("A19").Value = '4812574'
$ws.Range("B19").Value = 'https://vnexpress.net/thu-tuong-israel-sa-thai-bo-truong-quoc-phong-4812574.html'
$ws.Range("C19").Value = 'Thủ tướng Israel sa thải Bộ trưởng Quốc phòng'
$ws.Range("A20").Value = '4812524'
$ws.Range("B20").Value = 'https://vnexpress.net/nguoi-dan-khap-the-gioi-theo-doi-bau-cu-tong-thong-my-4812524.html'
$ws.Range("C20").Value = 'Người dân khắp thế giới theo dõi bầu cử tổng thống Mỹ'
$ws.Range("A21").Value = '4812566'
$ws.Range("B21").Value = 'https://vnexpress.net/nga-cao-buoc-ukraine-dinh-chiem-nha-may-dien-hat-nhan-o-kursk-4812566.html'
$ws.Range("C21").Value = 'Nga cáo buộc Ukraine định chiếm nhà máy điện hạt nhân ở Kursk'
$ws.Range("A22").Value = '4812516'
$ws.Range("B22").Value = 'https://vnexpress.net/noi-stress-cua-cu-tri-my-trong-ky-bau-cu-4812516.html'
$ws.Range("C22").Value = 'Nỗi stress của cử tri Mỹ trong kỳ bầu cử'
$ws.Range("A23").Value = '4812442'
$ws.Range("B23").Value = 'https://vnexpress.net/cach-gioi-chuc-my-ngan-phieu-bau-ma-4812442.html'
$ws.Range("C23").Value = 'Cách giới chức Mỹ ngăn ''phiếu bầu ma'''
$ws.Range("A24").Value = '4810347'
$ws.Range("B24").Value = 'https://vnexpress.net/ly-do-my-kho-cong-bo-nguoi-chien-thang-ngay-dem-bau-cu-4810347.html'
$ws.Range("C24").Value = 'Lý do Mỹ khó công bố người chiến thắng ngay đêm bầu cử'
$ws.Range("A25").Value = '4812551'
$ws.Range("B25").Value = 'https://vnexpress.net/ly-do-ong-trump-van-duoc-bo-phieu-du-da-bi-ket-toi-4812551.html'
$ws.Range("C25").Value = 'Lý do ông Trump vẫn được bỏ phiếu dù đã bị kết tội'
$ws.Range("A26").Value = '4812547'
$ws.Range("B26").Value = 'https://vnexpress.net/tiem-kich-an-do-xoay-nhu-la-vang-khi-roi-xuong-dat-4812547.html'
$ws.Range("C26").Value = 'Tiêm kích Ấn Độ ''xoay như lá vàng'' khi rơi xuống đất'
$ws.Range("A27").Value = '4812478'
$ws.Range("B27").Value = 'https://vnexpress.net/nguoi-an-do-toi-den-tho-cau-nguyen-cho-ba-harris-dac-cu-4812478.html'
$ws.Range("C27").Value = 'Người Ấn Độ tới đền thờ cầu nguyện cho bà Harris đắc cử'
$ws.Range("A28").Value = '4812154'
$ws.Range("B28").Value = 'https://vnexpress.net/an-so-tu-cu-tri-tham-lang-trong-bau-cu-my-4812154.html'
$ws.Range("C28").Value = 'Ẩn số từ cử tri thầm lặng trong bầu cử Mỹ'
$ws.Range("A29").Value = '4811827'
$ws.Range("B29").Value = 'https://vnexpress.net/may-bo-phieu-bau-tong-thong-my-hoat-dong-nhu-the-nao-4811827.html'
$ws.Range("C29").Value = 'Máy bỏ phiếu bầu tổng thống Mỹ hoạt động như thế nào?'
$ws.Range("A30").Value = '4812341'
$ws.Range("B30").Value = 'https://vnexpress.net/phat-bieu-khep-lai-chien-dich-tranh-cu-cua-trump-harris-4812341.html'
$ws.Range("C30").Value = 'Phát biểu khép lại chiến dịch tranh cử của Trump - Harris'
$ws.Range("A31").Value = '4812512'
$ws.Range("B31").Value = 'https://vnexpress.net/thu-tuong-tang-chan-dung-chu-tich-ho-chi-minh-cho-khu-di-tich-o-van-nam-4812512.html'
$ws.Range("C31").Value = 'Thủ tướng tặng chân dung Chủ tịch Hồ Chí Minh cho khu di tích ở Vân Nam'
$ws.Range("A32").Value = '4812455'
$ws.Range("B32").Value = 'https://vnexpress.net/cac-diem-bo-phieu-bau-tong-thong-tren-khap-nuoc-my-mo-cua-4812455.html'
$ws.Range("C32").Value = 'Ông Trump cáo buộc ''gian lận tràn lan'' ở Philadelphia'
$ws.Range("A33").Value = '4812468'
$ws.Range("B33").Value = 'https://vnexpress.net/vua-tay-ban-nha-duoc-ca-ngoi-vi-dung-vung-truoc-dam-dong-nem-bun-4812468.html'
$ws.Range("C33").Value = 'Vua Tây Ban Nha được ca ngợi vì đứng vững trước đám đông ném bùn'
$ws.Range("A34").Value = '4812291'
$ws.Range("B34").Value = 'https://vnexpress.net/ong-zelensky-nga-phong-uav-vao-ukraine-nhieu-gap-10-lan-nam-ngoai-4812291.html'
$ws.Range("C34").Value = 'Ông Zelensky: Nga phóng UAV vào Ukraine nhiều gấp 10 lần năm ngoái'
$ws.Range("A35").Value = '4812405'
$ws.Range("B35").Value = 'https://vnexpress.net/ha-ma-ngoi-sao-cua-thai-lan-du-doan-ong-trump-dac-cu-4812405.html'
$ws.Range("C35").Value = 'Hà mã ''ngôi sao'' của Thái Lan dự đoán ông Trump đắc cử'
$ws.Range("A36").Value = '4812390'
$ws.Range("B36").Value = 'https://vnexpress.net/ong-trump-chon-thanh-pho-dac-biet-lam-diem-van-dong-cuoi-cung-4812390.html'
$ws.Range("C36").Value = 'Ông Trump chọn ''thành phố đặc biệt'' làm điểm vận động cuối cùng'
$ws.Range("A37").Value = '4812160'
$ws.Range("B37").Value = 'https://vnexpress.net/nhung-gia-dinh-my-bat-dong-quan-diem-vi-bau-cu-tong-thong-4812160.html'
$ws.Range("C37").Value = 'Những gia đình Mỹ bất đồng quan điểm vì bầu cử tổng thống'
$ws.Range("A38").Value = '4812167'
$ws.Range("B38").Value = 'https://vnexpress.net/uav-lancet-co-the-da-tap-kich-hon-2-500-muc-tieu-o-ukraine-4812167.html'
$ws.Range("C38").Value = 'UAV Lancet có thể đã tập kích hơn 2.500 mục tiêu ở Ukraine'
$ws.Range("A39").Value = '4812177'
$ws.Range("B39").Value = 'https://vnexpress.net/nhung-hinh-anh-dinh-hinh-mua-bau-cu-tong-thong-my-4812177.html'
$ws.Range("C39").Value = 'Những hình ảnh định hình mùa bầu cử tổng thống Mỹ'
$ws.Range("A40").Value = '4812264'
$ws.Range("B40").Value = 'https://vnexpress.net/nguoi-dan-chuong-trinh-podcast-noi-tieng-tuyen-bo-ung-ho-ong-trump-4812264.html'
$ws.Range("C40").Value = 'Người dẫn chương trình podcast nổi tiếng tuyên bố ủng hộ ông Trump'
$ws.Range("A41").Value = '4812201'
$ws.Range("B41").Value = 'https://vnexpress.net/ong-trump-noi-ba-harris-nen-dau-voi-mike-tyson-4812201.html'
$ws.Range("C41").Value = 'Ông Trump nói bà Harris ''nên đấu với Mike Tyson'''
$ws.Range("A42").Value = '4811297'
$ws.Range("B42").Value = 'https://vnexpress.net/cach-nguoi-my-bo-phieu-va-kiem-dem-ket-qua-bau-cu-4811297.html'
$ws.Range("C42").Value = 'Cách người Mỹ bỏ phiếu và kiểm đếm kết quả bầu cử'
